$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Manchester City
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Manchester City"
$ws.Range("C2").Value = "England"
$ws.Range("D2").Value = "Premier Legue"
$ws.Range("E2").Value = "/uploads/clubs/1766017335291-earth-23593_1280.png"
$ws.Range("F2").Value = "man@gmail.com"
$ws.Range("G2").Value = "+11 77 89834234"
$ws.Range("H2").Value = "manchester.com"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0

# Row 3 - War Men
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "War Men"
$ws.Range("C3").Value = "Freetown"
$ws.Range("D3").Value = "Sierra Leone Premier"
$ws.Range("E3").Value = "/uploads/clubs/1766018020521-Tom___Jerry.jpeg"
$ws.Range("F3").Value = "warm@gmail.com"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "+323299888777"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = "warm.com"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
